# Add a new slide ("Title and Content" layout) after the existing slide,
# recreating the "Presentation" slide with the screen-capture link note.

$p = $ppt.ActivePresentation

# Index 2 -> insert after slide 1; layout 2 == ppLayoutText ("Title and Content").
$s = $p.Slides.Add(2, 2)

# --- Title -------------------------------------------------------------
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Text = "Presentation"

# --- Body copy -----------------------------------------------------------
$body = $s.Shapes.Item(2).TextFrame.TextRange

$part1 = "My Computer would not let me input my screen capture video into this "
$part2 = "powerpoint"
$part3 = ", but I was able to create this link to a location that would play the video."
$lead = "  "
$url = "https://docs.google.com/file/d/1JGqijwCs53T9afERzwA9crRNNYYaz7gz/preview"
$trail = "  "

# Paragraph 1: three separate runs.
$body.Text = $part1
$run2 = $body.InsertAfter($part2)
$run3 = $run2.InsertAfter($part3)

# Paragraph break, then paragraph 2: leading spaces / link / trailing spaces.
$run3.InsertAfter([char]13) | Out-Null
$runLead = $body.InsertAfter($lead)
$runUrl = $runLead.InsertAfter($url)
$runTrail = $runUrl.InsertAfter($trail)

# Apply the hyperlink to just the URL run (Characters() scopes correctly,
# unlike ActionSettings on an InsertAfter-chained range).
$linkStart = $part1.Length + $part2.Length + $part3.Length + 1 + $lead.Length + 1
$full = $s.Shapes.Item(2).TextFrame.TextRange
$link = $full.Characters($linkStart, $url.Length)
$link.ActionSettings.Item(1).Hyperlink.Address = $url

# --- Three leftover rectangle shapes (blank, space-only) ----------------
for ($i = 0; $i -lt 3; $i++) {
    $rect = $s.Shapes.AddShape(1, 349.5110236220472, 255.45937007874016, 20.977952755905513, 29.081259842519685)
    $rect.TextFrame.TextRange.Text = " "
    $rect.TextFrame.WordWrap = 0
    $rect.TextFrame.AutoSize = 1
}
